$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns.
# Numeric-looking price strings are prefixed with a literal quote so Excel
# keeps them as text (preserving trailing zeros / exact formatting) instead
# of auto-converting them to numbers.

$ws.Range("D2").Value = "43.108.70"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.369.94"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'304.14"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.482"
$ws.Range("E9").Value = "  -3.04%  "

$ws.Range("D10").Value = "'34.40"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("E13").Value = "  -2.32%  "

$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").Value = "2.736.81"
$ws.Range("E15").Value = "  +1.63%  "

$ws.Range("D16").Value = "2.359.92"
$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").Value = "43.104.86"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "'11.98"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").Value = "'68.15"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'235.46"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  -2.74%  "

$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'24.61"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("D30").Value = "'32.40"
$ws.Range("E30").Value = "  +2.43%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "'17.53"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").Value = "'0.0730"
$ws.Range("E34").Value = "  +4.05%  "

$ws.Range("E35").Value = "  +6.37%  "

$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("D37").Value = "'126.34"
$ws.Range("E37").Value = "  -10.44%  "

$ws.Range("D38").Value = "'4.34"
$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").Value = "'20.84"
$ws.Range("E42").Value = "  -6.93%  "

$ws.Range("D43").Value = "1.936.76"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").Value = "'2.15"
$ws.Range("E45").Value = "  +4.12%  "

$ws.Range("D46").Value = "'9.30"
$ws.Range("E46").Value = "  -7.80%  "

$ws.Range("D47").Value = "'2.72"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").Value = "2.596.13"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("E49").Value = "  +2.64%  "

$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").Value = "'1.14"
$ws.Range("E51").Value = "  +1.10%  "
